$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClueLayout")

# Swap room-label marker cells (room label "#" tiles moved to a new cell within each room)
$ws.Range("J4").Value = "H#"

$ws.Range("U11").Value = "P#"
$ws.Range("V11").Value = "P"

$ws.Range("C12").Value = "M#"
$ws.Range("O12").Value = "G#"
$ws.Range("Q12").Value = "G"
$ws.Range("D13").Value = "M"

$ws.Range("U18").Value = "S#"
$ws.Range("V18").Value = "S"
$ws.Range("U19").Value = "S*"

$ws.Range("C23").Value = "K#"
$ws.Range("L23").Value = "D#"
$ws.Range("P23").Value = "D"
$ws.Range("D24").Value = "K"

# Update the active selection on the sheet
$ws.Range("O21").Select()
